$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in A1 from "GRUPOS12" to "Valor"
$ws.Range("A1").Value = "Valor"

# Center-align headers A1:B1
$ws.Range("A1:B1").HorizontalAlignment = -4108  # xlCenter

# Set the active selection to C11 (as reflected in saved view state)
$ws.Range("C11").Select()
